$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename rule r16.21 from "Ammo Amounts" to "Ammo Limits" (cell B22).
#    Edited first so its shared-string entry is appended before the r4.3 one.
# ---------------------------------------------------------------------------
$newAmmoLimits = @'
<Bold>r16.21 Ammo Limits</Bold> 
<LineBreak/><LineBreak/>
Ammo is loaded during the Morning Briefing phase 
<InlineUIContainer><Button Content='r4.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  phase showing on the After Action Report 
<InlineUIContainer><Button Content='r2.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
how much of each ammo type is carried. 
<LineBreak/><LineBreak/>
Consult the 
<InlineUIContainer><Button Content='Ammo' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table for the amount of each type is available.
'@

$ws.Range("B22").Value2 = $newAmmoLimits

# ---------------------------------------------------------------------------
# 2) Add a 'Time' table button reference to rule r4.3 (cell B14).
# ---------------------------------------------------------------------------
$newTimeCheck = @'
<Bold>r4.3 Time Check Using the Time Tables</Bold> <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
Determine the hours of sunrise and sunset for the current month. On the After Action Report (AAR), the time is blocked out accordingly. Only the remaining time is available for play this day.
<LineBreak/><LineBreak/>
Determine the time and ammo already expended by rolling 1D on the Time Elapsed Table. This roll determines the number of hours already passed and the 
amount of ammo already expected before the scenario begins. The AAR shows the elapsed time and the ammo used up. 
'@

$ws.Range("B14").Value2 = $newTimeCheck

# ---------------------------------------------------------------------------
# 3) Update the view state: scroll so A11 is the top-left cell and select B15.
# ---------------------------------------------------------------------------
$ws.Range("B15").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
